$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure numeric-looking price strings (column D) stay text, matching the source data
# (which stores prices/percentages as literal strings, not numbers/dates).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.419.97'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.868.14'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.99'
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7074'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07874'
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3135'
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.56'
$ws.Range("E10").Value = '  -1.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07964'
$ws.Range("E11").Value = '  -4.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.906.55'
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.211'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.45'
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7013'
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.527'
$ws.Range("E16").Value = '  +2.65%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.464.84'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008362'
$ws.Range("E18").Value = '  -1.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.71'
$ws.Range("E19").Value = '  +3.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.125.28'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.13'
$ws.Range("E21").Value = '  -1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.640'
$ws.Range("E23").Value = '  -1.94%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1555'
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.004'
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.50'
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.72'
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.501'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.332'
$ws.Range("E30").Value = '  -1.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.257'
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05319'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.893'
$ws.Range("E34").Value = '  -2.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7498'
$ws.Range("E35").Value = '  -2.22%  '
$ws.Range("E36").Value = '  -1.00%  '
$ws.Range("E37").Value = '  +0.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01892'
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.273.90'
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.751'
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8945'
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.078'
$ws.Range("E42").Value = '  -6.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '109.18'
$ws.Range("E43").Value = '  -3.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.36'
$ws.Range("E44").Value = '  -4.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("E46").Value = '  -3.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.024.40'
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.797'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.563'
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.5182'
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4309'
$ws.Range("E51").Value = '  -1.44%  '
